$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $rng = $ws.Range($addr)
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = $origStyle
}

Set-TextValue "D2" '67.012.46'
Set-TextValue "E2" '  -0.50%  '
Set-TextValue "D3" '2.609.15'
Set-TextValue "E3" '  -1.30%  '
Set-TextValue "E4" '  -0.01%  '
Set-TextValue "D5" '590.62'
Set-TextValue "E5" '  -1.42%  '
Set-TextValue "D6" '165.42'
Set-TextValue "E6" '  -0.33%  '
Set-TextValue "E7" '  +0.01%  '
Set-TextValue "E8" '  -2.11%  '
Set-TextValue "D9" '2.608.43'
Set-TextValue "E9" '  -1.33%  '
Set-TextValue "E10" '  -4.86%  '
Set-TextValue "D11" '0.161'
Set-TextValue "E11" '  +1.58%  '
Set-TextValue "D12" '0.362'
Set-TextValue "E12" '  -0.50%  '
Set-TextValue "D13" '5.19'
Set-TextValue "E13" '  -0.71%  '
Set-TextValue "D14" '27.27'
Set-TextValue "E14" '  -2.53%  '
Set-TextValue "D15" '3.084.29'
Set-TextValue "E15" '  -1.28%  '
Set-TextValue "E16" '  -2.62%  '
Set-TextValue "D17" '66.945.37'
Set-TextValue "E17" '  -0.39%  '
Set-TextValue "D18" '2.630.98'
Set-TextValue "E18" '  -0.50%  '
Set-TextValue "D19" '11.78'
Set-TextValue "E19" '  -0.95%  '
Set-TextValue "D20" '7.80'
Set-TextValue "E20" '  -0.92%  '
Set-TextValue "D21" '354.17'
Set-TextValue "E21" '  -2.68%  '
Set-TextValue "D22" '4.27'
Set-TextValue "E22" '  -2.90%  '
Set-TextValue "D23" '4.62'
Set-TextValue "E23" '  -3.47%  '
Set-TextValue "D24" '10.53'
Set-TextValue "E24" '  -5.12%  '
Set-TextValue "E25" '  +0.01%  '
Set-TextValue "E26" '  -4.66%  '
Set-TextValue "D27" '68.97'
Set-TextValue "E27" '  -2.68%  '
Set-TextValue "E28" '  -1.06%  '
Set-TextValue "E29" '  -0.01%  '
Set-TextValue "D30" '0.0₃0995'
Set-TextValue "E30" '  -2.83%  '
Set-TextValue "D31" '540.43'
Set-TextValue "E31" '  -2.49%  '
Set-TextValue "D32" '7.86'
Set-TextValue "E32" '  -2.22%  '
Set-TextValue "E34" '  -2.85%  '
Set-TextValue "D35" '0.134'
Set-TextValue "E35" '  +0.70%  '
Set-TextValue "E36" '  -0.01%  '
Set-TextValue "D37" '1.49'
Set-TextValue "E37" '  -3.69%  '
Set-TextValue "D38" '157.16'
Set-TextValue "E38" '  -0.33%  '
Set-TextValue "D39" '18.90'
Set-TextValue "E39" '  -2.65%  '
Set-TextValue "E40" '  -2.21%  '
Set-TextValue "D41" '18.23'
Set-TextValue "E41" '  +1.79%  '
Set-TextValue "E42" '  -1.20%  '
Set-TextValue "D43" '5.13'
Set-TextValue "E43" '  -2.57%  '
Set-TextValue "E45" '  -4.76%  '
Set-TextValue "E46" '  -1.52%  '
Set-TextValue "D47" '151.07'
Set-TextValue "E47" '  -1.97%  '
Set-TextValue "D48" '0.574'
Set-TextValue "E48" '  -3.72%  '
Set-TextValue "E49" '  -3.24%  '
Set-TextValue "D50" '1.70'
Set-TextValue "E50" '  -2.22%  '
Set-TextValue "E51" '  -1.14%  '
